$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rf-adjusted statistics (B, D, E, F) for rows 2-6
# Row 2: HML
$ws.Range("B2").Value = 12.71976515605532
$ws.Range("D2").Value = 0.009947741319881393
$ws.Range("E2").Value = 1.463051787562091
$ws.Range("F2").Value = 2.140520533088629

# Row 3: HML RMW
$ws.Range("B3").Value = 12.74248969159571
$ws.Range("D3").Value = 0.009936579308624488
$ws.Range("E3").Value = 1.461410148520918
$ws.Range("F3").Value = 2.135719622199932

# Row 4: HML CMA
$ws.Range("B4").Value = 13.33577442481397
$ws.Range("D4").Value = 0.009445451835545612
$ws.Range("E4").Value = 1.389178181051797
$ws.Range("F4").Value = 1.92981601871038

# Row 5: RMW CMA
$ws.Range("B5").Value = 13.34703001737696
$ws.Range("D5").Value = 0.009888440638458064
$ws.Range("E5").Value = 1.454330212968403
$ws.Range("F5").Value = 2.115076368352721

# Row 6: HML RMW CMA
$ws.Range("B6").Value = 13.34209866551272
$ws.Range("D6").Value = 0.00943748391312205
$ws.Range("E6").Value = 1.388006308686997
$ws.Range("F6").Value = 1.926561512954904
